# Update odds values in Sheet1 per the 2025-01-29 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.45
$ws.Range("H2").Value = 2.8
$ws.Range("O2").Value = 1.67
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 2.39
$ws.Range("R2").Value = 1.58
$ws.Range("U2").Value = 5.8
$ws.Range("V2").Value = 1.14
$ws.Range("W2").Value = 7
$ws.Range("X2").Value = 1.1

# Row 3
$ws.Range("G3").Value = 2.25
$ws.Range("H3").Value = 2.9

# Row 4
$ws.Range("I4").Value = 2.4

# Row 5
$ws.Range("G5").Value = 2.9
$ws.Range("I5").Value = 2.75
$ws.Range("J5").Value = 3.75
$ws.Range("M5").Value = 1.14
$ws.Range("N5").Value = 5.5
$ws.Range("Y5").Value = 1.67
$ws.Range("Z5").Value = 2.1
$ws.Range("AD5").Value = 12
$ws.Range("AF5").Value = 29
$ws.Range("AG5").Value = 29
$ws.Range("AO5").Value = 12
$ws.Range("AP5").Value = 29

# Row 8
$ws.Range("S8").Value = 2
$ws.Range("T8").Value = 1.85
$ws.Range("W8").Value = 3.4
$ws.Range("X8").Value = 1.3
$ws.Range("Y8").Value = 1.37

# Row 13
$ws.Range("AJ13").Value = 6

# Row 17
$ws.Range("G17").Value = 2.38
$ws.Range("H17").Value = 2.88
$ws.Range("I17").Value = 3.5
$ws.Range("J17").Value = 3.25
$ws.Range("M17").Value = 1.13
$ws.Range("N17").Value = 6
$ws.Range("AC17").Value = 6
$ws.Range("AD17").Value = 10
$ws.Range("AE17").Value = 10
$ws.Range("AM17").Value = 8
$ws.Range("AQ17").Value = 34
